$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.026.57"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.479.19"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.19%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "319.59"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "107.63"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +2.38%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.520"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  +0.06%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.530"
$r.Style = "Normal"
$ws.Range("E9").Value = "  -1.28%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "38.38"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +6.40%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0806"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("E12").Value = "  +0.21%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "18.08"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -1.21%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "7.08"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "2.874.98"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "2.488.69"
$ws.Range("E16").Value = "  +1.41%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.840"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "46.955.13"
$ws.Range("E18").Value = "  +1.01%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "12.62"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "6.56"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("E21").Value = "  +15.57%  "
$ws.Range("D22").Value = "0.0₃0929"
$ws.Range("E22").Value = "  -0.73%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "70.47"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "244.47"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -1.82%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "2.54"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").Value = "  +0.01%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "25.51"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "2.24"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -3.03%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "9.98"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("E30").Value = "  +8.24%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "34.69"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -1.37%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "49.66"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "19.83"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +1.09%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "5.31"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.0778"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  +2.14%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "4.60"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -0.52%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "2.93"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.111"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  -0.31%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "119.10"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("E43").Value = "  +3.22%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.0292"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "1.972.91"
$ws.Range("E45").Value = "  -0.45%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "2.99"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.53%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "1.99"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -3.83%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "9.01"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +0.83%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.74"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -3.26%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "5.10"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -4.71%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "56.53"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +3.95%  "
